$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that is bumped by one day
# (from 46074 -> 46075) for every data row (rows 2 through 173).
$ws.Range("C2:C173").Value = 46075
